# Applies the "Boll, brickor, kollision, pausmeny" edit to the test spec.
$d = $word.ActiveDocument

function Replace-ParaXml($para, [string]$xml) {
    $r = $para.Range
    $r.InsertXML($xml) | Out-Null
}

# 1) Title: merge "Testspecifikation" + " ”Break a Brick”" into a single run.
$p1 = $d.Paragraphs.Item(1)
Replace-ParaXml $p1 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Testspecifikation ”Break a Brick”</w:t></w:r></w:p>'

# 2) TF 1.1/1.2 heading: merge "TF 1" + ".2" + " " into one run attached to the line break.
$p7 = $d.Paragraphs.Item(7)
Replace-ParaXml $p7 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>TF 1.1 Meny:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Huvudmeny</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">TF 1.2 </w:t></w:r><w:r><w:t>Meny: Starta spelet</w:t></w:r></w:p>'

# 3) "Systemtest Spel" heading: merge "Systemtest " + "Spel" into a single run.
$p12 = $d.Paragraphs.Item(12)
Replace-ParaXml $p12 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Systemtest Spel</w:t></w:r></w:p>'

# 4) Renumber "TF 2.1 Spel: Plattan" -> "TF 2 Spelet" and insert the whole new
#    "Systemtest Pausmeny" section (heading, intro, Testfall, "TF 3 Pausmeny")
#    plus the extra trailing blank paragraphs, replacing the old paragraphs
#    16-19 (TF2 heading .. the two blanks .. the blank Heading4 paragraph).
$p16 = $d.Paragraphs.Item(16)
$p19 = $d.Paragraphs.Item(19)
$r16 = $d.Range($p16.Range.Start, $p19.Range.End)
$r16.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>TF 2 Spelet</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">Systemtest </w:t></w:r><w:r><w:t>Pausmeny</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Systemtester ska testa så att </w:t></w:r><w:r><w:t xml:space="preserve">pausmenyn </w:t></w:r><w:r><w:t>fungerar på tänkt sätt.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Testfall</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>TF 3 Pausmeny</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" /><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr></w:p>') | Out-Null
